$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# row, nombre_aides (C), nombre_entreprises (D), montant_total (E)
$updates = @(
    @{Row=20;  C=108353; D=24272; E=225644598},
    @{Row=23;  C=12431;  D=3244;  E=26551383},
    @{Row=97;  C=98506;  D=22696; E=307112902},
    @{Row=164; C=50570;  D=11057; E=168509301},
    @{Row=168; C=284949; D=58122; E=1208950823},
    @{Row=169; C=562571; D=60954; E=1284349419},
    @{Row=170; C=367306; D=38109; E=2844387745},
    @{Row=171; C=115122; D=20263; E=445187670},
    @{Row=173; C=54385;  D=11601; E=151851439},
    @{Row=174; C=357180; D=69788; E=1016827280},
    @{Row=175; C=125519; D=18100; E=812155990},
    @{Row=177; C=96749;  D=16507; E=174720542},
    @{Row=179; C=235668; D=29335; E=812410531},
    @{Row=180; C=141464; D=28890; E=340163580},
    @{Row=193; C=5346;   D=1046;  E=27703233},
    @{Row=210; C=6419;   D=979;   E=18553407},
    @{Row=213; C=3635;   D=402;   E=11266871},
    @{Row=266; C=71664;  D=9133;  E=219431842},
    @{Row=322; C=81162;  D=9703;  E=254537631}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
